$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (fSorb) values with new stock calculation results
$ws.Range("B2").Value = [double]"-5.8477646057199999E-4"
$ws.Range("B3").Value = [double]"-5.8477142137400003E-4"
$ws.Range("B4").Value = [double]"-5.8476679864299995E-4"
$ws.Range("B5").Value = [double]"-5.8476063208000003E-4"
$ws.Range("B6").Value = [double]"-5.8474932247400005E-4"
$ws.Range("B7").Value = [double]"-5.8472605153099999E-4"
$ws.Range("B8").Value = [double]"-5.8467810525300004E-4"
$ws.Range("B9").Value = [double]"-5.8458251305200002E-4"
$ws.Range("B10").Value = [double]"-5.8439916722700003E-4"
$ws.Range("B11").Value = [double]"-5.8406019795899996E-4"
$ws.Range("B12").Value = [double]"-5.8345445943599995E-4"
$ws.Range("B13").Value = [double]"-5.8240100784400005E-4"
$ws.Range("B14").Value = [double]"-5.8061776242600001E-4"
$ws.Range("B15").Value = [double]"-5.7766749215699997E-4"
$ws.Range("B16").Value = [double]"-5.7288529583599999E-4"
$ws.Range("B17").Value = [double]"-5.6527500260199995E-4"
$ws.Range("B18").Value = [double]"-5.5336448200599998E-4"
$ws.Range("B19").Value = [double]"-5.3500606705400003E-4"
$ws.Range("B20").Value = [double]"-5.0710485226299995E-4"
$ws.Range("B21").Value = [double]"-4.6525402617699998E-4"
$ws.Range("B22").Value = [double]"-4.0325316161000002E-4"
$ws.Range("B23").Value = [double]"-3.1248347782600002E-4"
$ws.Range("B24").Value = [double]"-1.8111509638000001E-4"
$ws.Range("B25").Value = [double]"6.8723226201199997E-6"
$ws.Range("B26").Value = [double]"2.7285898286200002E-4"
$ws.Range("B27").Value = [double]"6.4491927966899996E-4"
$ws.Range("B28").Value = [double]"1.15924894859E-3"
$ws.Range("B29").Value = [double]"1.8615429855999999E-3"
$ws.Range("B30").Value = [double]"2.8080881075E-3"
$ws.Range("B31").Value = [double]"4.0662213893400004E-3"
$ws.Range("B32").Value = [double]"5.7136965936499997E-3"
$ws.Range("B33").Value = [double]"7.8364316155299998E-3"
$ws.Range("B34").Value = [double]"1.0524147136799999E-2"
$ws.Range("B35").Value = [double]"1.3863625862699999E-2"
$ws.Range("B36").Value = [double]"1.7929788093900001E-2"
$ws.Range("B37").Value = [double]"2.2775496357800001E-2"
$ws.Range("B38").Value = [double]"2.8421847812099998E-2"
$ws.Range("B39").Value = [double]"3.4851412099699997E-2"
$ws.Range("B40").Value = [double]"4.2007041730500003E-2"
$ws.Range("B41").Value = [double]"4.9798200749299998E-2"
$ws.Range("B42").Value = [double]"5.8115198723E-2"
$ws.Range("B43").Value = [double]"6.6849697721499998E-2"
$ws.Range("B44").Value = [double]"7.5918138324000004E-2"
$ws.Range("B45").Value = [double]"8.5284036308499997E-2"
$ws.Range("B46").Value = [double]"9.4975710433500005E-2"
$ws.Range("B47").Value = [double]"0.105097564273"
$ws.Range("B48").Value = [double]"0.115834801369"
$ws.Range("B49").Value = [double]"0.127452664203"
$ws.Range("B50").Value = [double]"0.14029157624999999"
$ws.Range("B51").Value = [double]"0.154758984307"
$ws.Range("B52").Value = [double]"0.17131757577500001"
$ws.Range("B53").Value = [double]"0.19046830073500001"
$ws.Range("B54").Value = [double]"0.21272566936599999"
$ws.Range("B55").Value = [double]"0.23858252425900001"
$ws.Range("B56").Value = [double]"0.26846234310900002"
$ws.Range("B57").Value = [double]"0.30265952501400001"
$ws.Range("B58").Value = [double]"0.34127222674000002"
$ws.Range("B59").Value = [double]"0.384137591793"
$ws.Range("B60").Value = [double]"0.430784590099"
$ws.Range("B61").Value = [double]"0.48042065554199997"
$ws.Range("B62").Value = [double]"0.53196644757900002"
$ws.Range("B63").Value = [double]"0.58414121044099998"
$ws.Range("B64").Value = [double]"0.63558776122899996"
$ws.Range("B65").Value = [double]"0.68501290854800001"
$ws.Range("B66").Value = [double]"0.73131377720599999"
$ws.Range("B67").Value = [double]"0.77366544613199995"
$ws.Range("B68").Value = [double]"0.81155802083999995"
$ws.Range("B69").Value = [double]"0.844785735752"
$ws.Range("B70").Value = [double]"0.87340116331999995"
$ws.Range("B71").Value = [double]"0.89765150550999995"
$ws.Range("B72").Value = [double]"0.91791211648500004"
$ws.Range("B73").Value = [double]"0.93462753358100004"
$ws.Range("B74").Value = [double]"0.948265015268"
$ws.Range("B75").Value = [double]"0.95928149768600002"
$ws.Range("B76").Value = [double]"0.96810245692300001"
$ws.Range("B77").Value = [double]"0.97511016721099997"
$ws.Range("B78").Value = [double]"0.98063877047000003"
$ws.Range("B79").Value = [double]"0.98497394685100004"
$ws.Range("B80").Value = [double]"0.98835547923400002"
$ws.Range("B81").Value = [double]"0.99098147154399996"
$ws.Range("B82").Value = [double]"0.99301335446799999"

# Update column A (pH) values that changed precision
$ws.Range("A62").Value = [double]"7.9999999999"
$ws.Range("A73").Value = [double]"9.0999999997999996"

# Update the active selection cell
$ws.Range("D55").Select()
